$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,8).Value = "land"
